$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows above the current row 2, pushing the existing
# rows 2-15 down to rows 9-22. Then strip any formatting Excel may
# have copied onto the freshly-inserted rows so they start out
# "un-styled", matching the rest of the data rows.
$ws.Rows("2:8").Insert()
$ws.Range("A2:E8").ClearFormats()

# New processes that belong at the top of the list.
$newTopRows = @(
    @{ Row = 2;  A = "2024/112521.7"; B = "SECRETARIA 2-B"; C = 45331; D = "SANDERLAN"; E = "Fulano5" },
    @{ Row = 3;  A = "2024/113511.8"; B = "SECRETARIA 2-B"; C = 45355; D = "MARCOS";    E = "Fulano2" },
    @{ Row = 4;  A = "2024/113381.2"; B = "SECRETARIA 2-B"; C = 45355; D = "SABRINA";   E = "Fulano1`n" },
    @{ Row = 5;  A = "2024/112385.2"; B = "SECRETARIA 2-B"; C = 45579; D = "SANDERLAN"; E = "Fulano3" },
    @{ Row = 6;  A = "2024/413385.6"; B = "SECRETARIA 2-B"; C = 45576; D = "MARCOS";    E = "Fulano4`n" },
    @{ Row = 7;  A = "2024/553362.6"; B = "SECRETARIA 2-B"; C = 45576; D = "YGOR";      E = "Fulano2" },
    @{ Row = 8;  A = "2024/955322.1"; B = "SECRETARIA 2-B"; C = 45576; D = "EDUARDO";   E = "Fulano1`n" }
)

# New processes appended after the (now shifted) existing data, at the
# bottom of the sheet.
$newBottomRows = @(
    @{ Row = 23; A = "2024/015609.9"; B = "SECRETARIA 2-B"; C = 45391; D = "RUY";       E = "Fulano4`n" },
    @{ Row = 24; A = "2024/112521.8"; B = "SECRETARIA 2-B"; C = 45331; D = "SANDERLAN"; E = "Fulano3" }
)

foreach ($r in ($newTopRows + $newBottomRows)) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E

    # Cells whose "responsavel" value carries an embedded newline make
    # Excel auto-grow the row height; auto-fit it back down so the row
    # stays at the sheet's normal (non-custom) height, same as every
    # other data row in the sheet.
    $ws.Rows($row).EntireRow.AutoFit()
}
